# Append the new attendance ("absensi") log rows 17-30 to Sheet1.
# Source data scraped from the scanner log for 2024-09-21 (columns:
# Nama, NIS, Kelas, Waktu Kehadiran, Status).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column B holds NIS (student id) numbers, some of which have leading
# zeros (e.g. "0066944417", "0012", "00133"). Force that column to Text
# formatting before writing so Excel doesn't coerce the values to
# numbers and silently drop the leading zeros.
$ws.Range("B17:B30").NumberFormat = "@"

# Row 17
$ws.Range("B17").Value = "12225495"
$ws.Range("D17").Value = "14:20:11"
$ws.Range("E17").Value = "Terlambat"

# Row 18
$ws.Range("A18").Value = "surya sahrul"
$ws.Range("B18").Value = "3131"
$ws.Range("C18").Value = "12c1"
$ws.Range("D18").Value = "14:21:38"
$ws.Range("E18").Value = "Terlambat"

# Row 19
$ws.Range("A19").Value = "surya sahrul"
$ws.Range("B19").Value = "1222549"
$ws.Range("C19").Value = "12c1"
$ws.Range("D19").Value = "14:22:17"
$ws.Range("E19").Value = "Terlambat"

# Row 20
$ws.Range("A20").Value = "Tidak Diketahui"
$ws.Range("B20").Value = "31312"
$ws.Range("C20").Value = "12c1"
$ws.Range("D20").Value = "14:29:04"
$ws.Range("E20").Value = "Terlambat"

# Row 21
$ws.Range("A21").Value = "Tidak Diketahui"
$ws.Range("B21").Value = "313121"
$ws.Range("C21").Value = "12c1"
$ws.Range("D21").Value = "14:29:57"
$ws.Range("E21").Value = "Terlambat"

# Row 22
$ws.Range("B22").Value = "31312"
$ws.Range("D22").Value = "14:30:47"
$ws.Range("E22").Value = "Terlambat"

# Row 23
$ws.Range("A23").Value = "surya sahrul muhammad"
$ws.Range("B23").Value = "313121"
$ws.Range("C23").Value = "12c1"
$ws.Range("D23").Value = "14:31:39"
$ws.Range("E23").Value = "Terlambat"

# Row 24
$ws.Range("B24").Value = "122025495100"
$ws.Range("D24").Value = "14:38:50"
$ws.Range("E24").Value = "Terlambat"

# Row 25
$ws.Range("B25").Value = "122032545100"
$ws.Range("D25").Value = "14:38:50"
$ws.Range("E25").Value = "Terlambat"

# Row 26
$ws.Range("B26").Value = "12225495100"
$ws.Range("D26").Value = "14:38:52"
$ws.Range("E26").Value = "Terlambat"

# Row 27
$ws.Range("B27").Value = "0066944417"
$ws.Range("D27").Value = "14:38:56"
$ws.Range("E27").Value = "Terlambat"

# Row 28
$ws.Range("B28").Value = "1234"
$ws.Range("D28").Value = "14:38:58"
$ws.Range("E28").Value = "Terlambat"

# Row 29
$ws.Range("B29").Value = "0012"
$ws.Range("D29").Value = "14:47:52"
$ws.Range("E29").Value = "Terlambat"

# Row 30
$ws.Range("B30").Value = "00133"
$ws.Range("D30").Value = "14:48:17"
$ws.Range("E30").Value = "Terlambat"
